# TC10_Trials_Filter_Diagnosis-Colorectal.xlsx
# - added the Cypher/Neo4j query used to produce this test-case's data into A2
# - row 2 grows tall enough to show the (wrapped) query text
# - selection moves to A2 (the cell that now holds the query)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$query = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Colorectal cancer, NOS'] RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"

$ws.Range("A2").Value = $query

# Row 2 is tall enough now to display the wrapped query text.
$ws.Rows.Item(2).RowHeight = 87

# Selection follows the edited cell.
[void]$ws.Range("A2").Select()
